# Update MSME country indicator figures on the "Summary" sheet with more
# precise (extra decimal) values. The cells hold numeric-looking text
# (shared strings), so each cell is explicitly formatted as Text before
# the write to stop Excel from auto-coercing the string into a Number,
# then ClearFormats() restores the cell's original (General/default)
# style once the text value has been committed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

function Set-TextValue {
    param($ws, $addr, $text)
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

# Enterprises density (per 1000 people): Micro / SMEs
Set-TextValue $ws "B11" "10.23"
Set-TextValue $ws "C11" "4.17"

# Employment (% of total): Micro / SMEs
Set-TextValue $ws "B12" "11.53"
Set-TextValue $ws "C12" "39.58"

# Employment (absolute #): Micro / SMEs
Set-TextValue $ws "B13" "730308.25"
Set-TextValue $ws "C13" "2507698.75"

# Enterprises (% of total): SMEs / MSMEs
Set-TextValue $ws "C14" "28.44"
Set-TextValue $ws "D14" "98.14"
